$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# 1) Rows 10 and 11: the two observations got their Id (A), Ost (Q) and
#    Nord (R) values swapped with each other. Everything else in these
#    rows is identical between the two rows, so a direct value swap of
#    just these three columns is sufficient.
# -------------------------------------------------------------------------
$a10 = $ws.Range("A10").Value2
$a11 = $ws.Range("A11").Value2
$q10 = $ws.Range("Q10").Value2
$q11 = $ws.Range("Q11").Value2
$r10 = $ws.Range("R10").Value2
$r11 = $ws.Range("R11").Value2

$ws.Range("A10").Value = $a11
$ws.Range("A11").Value = $a10
$ws.Range("Q10").Value = $q11
$ws.Range("Q11").Value = $q10
$ws.Range("R10").Value = $r11
$ws.Range("R11").Value = $r10

# -------------------------------------------------------------------------
# 2) Rows 16 and 17: the two full observation records got swapped with
#    each other. Rather than copying whole rows (which risks turning the
#    text dates in Y/AA into real Excel date serials and disturbing
#    untouched columns), only the columns whose value actually differs
#    between the two rows are updated directly.
# -------------------------------------------------------------------------

# Columns with a plain value swap (non-blank on both sides)
$swapCols = @("A","B","E","F","G","H","P","Q","R","S","AW","AX")
foreach ($c in $swapCols) {
    $v16 = $ws.Range($c + "16").Value2
    $v17 = $ws.Range($c + "17").Value2
    $ws.Range($c + "16").Value = $v17
    $ws.Range($c + "17").Value = $v16
}

# M: row16 gains "aldre spar" (was blank), row17 loses it (becomes blank)
$ws.Range("M16").Value = $ws.Range("M17").Value2
$ws.Range("M17").ClearContents()

# AC: row16 gains "Ringhack pa gran" (was blank), row17 loses it
$ws.Range("AC16").Value = $ws.Range("AC17").Value2
$ws.Range("AC17").ClearContents()

# K, L, N: row16 gains empty placeholder cells (were entirely absent),
# row17's placeholder cells become entirely absent (were empty strings).
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("K17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("N17").ClearContents()

Write-Host "Row 10/11 swap and Row 16/17 swap complete."
